$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("openbis-data")

# Incorporate Bernd's changes that validate strains against the straindb:
# prefix the hybridization header labels with "JJS-" so they match the
# strain names used elsewhere (JJS-MGP253 / JJS-MGP776).
$ws.Range("D5").Value = "JJS-MGP253-1 66687802"
$ws.Range("E5").Value = "JJS-MGP776-2 66730002"

# Move the active selection on the data sheet.
$ws.Range("H11").Select()
